$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 839
$ws.Range("I28").Value = 435.17648
$ws.Range("K28").Value = 435.17648
$ws.Range("M28").Value = 49.82351999999997
$ws.Range("H32").Value = 16669783
$ws.Range("I32").Value = 1750
$ws.Range("J32").Value = 25003800
$ws.Range("K32").Value = 1750
$ws.Range("L32").Value = 25003800
$ws.Range("M32").Value = -1424
$ws.Range("N32").Value = -25004452
$ws.Range("H43").Value = 4214.5713
$ws.Range("J43").Value = 4200.4
$ws.Range("L43").Value = 4200.4
$ws.Range("N43").Value = -4338.4
$ws.Range("H80").Value = 2931.9355
$ws.Range("J80").Value = 3825.1365
$ws.Range("L80").Value = 11475.4095
$ws.Range("N80").Value = -13471.4095
$ws.Range("H83").Value = 2931.9355
$ws.Range("J83").Value = 3825.1365
$ws.Range("L83").Value = 34426.2285
$ws.Range("N83").Value = -44410.2285
$ws.Range("H98").Value = 1722.1482
$ws.Range("J98").Value = 938.6667
$ws.Range("L98").Value = 938.6667
$ws.Range("N98").Value = -3934.6667
$ws.Range("H100").Value = 8939.066000000001
$ws.Range("I100").Value = 7297.143
$ws.Range("J100").Value = 10375.75
$ws.Range("K100").Value = 7297.143
$ws.Range("L100").Value = 10375.75
$ws.Range("M100").Value = -6756.143
$ws.Range("N100").Value = -11457.75
$ws.Range("H101").Value = 666.5
$ws.Range("I101").Value = 668
$ws.Range("J101").Value = 663.5
$ws.Range("K101").Value = 2004
$ws.Range("L101").Value = 1990.5
$ws.Range("M101").Value = -382
$ws.Range("N101").Value = -5234.5
$ws.Range("H103").Value = 860.875
$ws.Range("I103").Value = 853
$ws.Range("J103").Value = 868.75
$ws.Range("K103").Value = 2559
$ws.Range("L103").Value = 2606.25
$ws.Range("M103").Value = -1973
$ws.Range("N103").Value = -3778.25
$ws.Range("H122").Value = 1722.1482
$ws.Range("J122").Value = 938.6667
$ws.Range("L122").Value = 2816.0001
$ws.Range("N122").Value = -7716.0001
$ws.Range("H138").Value = 4999.25
$ws.Range("I138").Value = 4998
$ws.Range("J138").Value = 4999.4287
$ws.Range("K138").Value = 14994
$ws.Range("L138").Value = 14998.2861
$ws.Range("M138").Value = -9854
$ws.Range("N138").Value = -25278.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5924.5586
$ws.Range("I32").Value = 6825.5615
$ws.Range("J32").Value = 1255.7273
$ws.Range("K32").Value = 6825.5615
$ws.Range("L32").Value = 1255.7273
$ws.Range("M32").Value = -6538.5615
$ws.Range("N32").Value = -1829.7273
$ws.Range("H45").Value = 4042.625
$ws.Range("I45").Value = 3191.6191
$ws.Range("J45").Value = 9999.666999999999
$ws.Range("K45").Value = 3191.6191
$ws.Range("L45").Value = 9999.666999999999
$ws.Range("M45").Value = -2814.6191
$ws.Range("N45").Value = -10753.667
$ws.Range("H132").Value = 1470.0454
$ws.Range("I132").Value = 1491.762
$ws.Range("J132").Value = 1014
$ws.Range("K132").Value = 4475.286
$ws.Range("L132").Value = 3042
$ws.Range("M132").Value = -1945.286
$ws.Range("N132").Value = -8102

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 252450
$ws.Range("J23").Value = 4900
$ws.Range("L23").Value = 4900
$ws.Range("N23").Value = -5466
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 144.04
$ws.Range("I7").Value = 78.882355
$ws.Range("J7").Value = 282.5
$ws.Range("K7").Value = 78.882355
$ws.Range("L7").Value = 282.5
$ws.Range("M7").Value = 34.117645
$ws.Range("N7").Value = -508.5
$ws.Range("H31").Value = 17007.428
$ws.Range("I31").Value = 34368.668
$ws.Range("J31").Value = 3986.5
$ws.Range("K31").Value = 34368.668
$ws.Range("L31").Value = 3986.5
$ws.Range("M31").Value = -34073.668
$ws.Range("N31").Value = -4576.5
$ws.Range("H34").Value = 17007.428
$ws.Range("I34").Value = 34368.668
$ws.Range("J34").Value = 3986.5
$ws.Range("K34").Value = 34368.668
$ws.Range("L34").Value = 3986.5
$ws.Range("M34").Value = -34166.668
$ws.Range("N34").Value = -4390.5
$ws.Range("H58").Value = 2183.3
$ws.Range("I58").Value = 1619.6666
$ws.Range("J58").Value = 2424.8572
$ws.Range("K58").Value = 1619.6666
$ws.Range("L58").Value = 2424.8572
$ws.Range("M58").Value = -1416.6666
$ws.Range("N58").Value = -2830.8572
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H99").Value = 14031360
$ws.Range("I99").Value = 2715861.8
$ws.Range("K99").Value = 2715861.8
$ws.Range("M99").Value = -2714363.8
$ws.Range("H126").Value = 14031360
$ws.Range("I126").Value = 2715861.8
$ws.Range("K126").Value = 8147585.399999999
$ws.Range("M126").Value = -8145115.399999999
$ws.Range("H136").Value = 2183.3
$ws.Range("I136").Value = 1619.6666
$ws.Range("J136").Value = 2424.8572
$ws.Range("K136").Value = 4858.9998
$ws.Range("L136").Value = 7274.571599999999
$ws.Range("M136").Value = -2308.9998
$ws.Range("N136").Value = -12374.5716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.636364
$ws.Range("I2").Value = 36.11111
$ws.Range("J2").Value = 99.5
$ws.Range("K2").Value = 216.66666
$ws.Range("L2").Value = 597
$ws.Range("M2").Value = -103.66666
$ws.Range("N2").Value = -823
$ws.Range("H5").Value = 630.3214
$ws.Range("I5").Value = 1100.3334
$ws.Range("K5").Value = 3301.0002
$ws.Range("M5").Value = -3189.0002
$ws.Range("H8").Value = 520.75
$ws.Range("I8").Value = 520.75
$ws.Range("K8").Value = 1562.25
$ws.Range("M8").Value = -1423.25
$ws.Range("H15").Value = 103.333336
$ws.Range("I15").Value = 150
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 450
$ws.Range("L15").Value = 30
$ws.Range("M15").Value = -310
$ws.Range("N15").Value = -310
$ws.Range("H40").Value = 17
$ws.Range("I40").Value = 6
$ws.Range("K40").Value = 24
$ws.Range("M40").Value = 45
$ws.Range("H61").Value = 466.66666
$ws.Range("H86").Value = 997.5
$ws.Range("I86").Value = 997.5
$ws.Range("K86").Value = 2992.5
$ws.Range("M86").Value = -1806.5
$ws.Range("H89").Value = 997.5
$ws.Range("I89").Value = 997.5
$ws.Range("K89").Value = 8977.5
$ws.Range("M89").Value = -3049.5
$ws.Range("H107").Value = 1689
$ws.Range("I107").Value = 1224.1538
$ws.Range("J107").Value = 2238.3635
$ws.Range("K107").Value = 3672.4614
$ws.Range("L107").Value = 6715.0905
$ws.Range("M107").Value = -1752.4614
$ws.Range("N107").Value = -10555.0905
$ws.Range("H113").Value = 672.93335
$ws.Range("I113").Value = 774.3333
$ws.Range("J113").Value = 520.8333
$ws.Range("K113").Value = 2322.9999
$ws.Range("L113").Value = 1562.4999
$ws.Range("M113").Value = -152.9998999999998
$ws.Range("N113").Value = -5902.4999
$ws.Range("H122").Value = 4269.5
$ws.Range("I122").Value = 892
$ws.Range("J122").Value = 5020.0557
$ws.Range("K122").Value = 8028
$ws.Range("L122").Value = 45180.5013
$ws.Range("M122").Value = -5578
$ws.Range("N122").Value = -50080.5013
$ws.Range("H132").Value = 1003.1905
$ws.Range("I132").Value = 1017.4667
$ws.Range("J132").Value = 967.5
$ws.Range("K132").Value = 9157.2003
$ws.Range("L132").Value = 8707.5
$ws.Range("M132").Value = -6627.2003
$ws.Range("N132").Value = -13767.5
$ws.Range("H135").Value = 630.3214
$ws.Range("I135").Value = 1100.3334
$ws.Range("K135").Value = 9903.000599999999
$ws.Range("M135").Value = -7368.000599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2531.6667
$ws.Range("I80").Value = 2438
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2438
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1440
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2531.6667
$ws.Range("I83").Value = 2438
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 12190
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -7198
$ws.Range("N83").Value = -24984
$ws.Range("H113").Value = 13673.523
$ws.Range("I113").Value = 10808.143
$ws.Range("J113").Value = 15106.214
$ws.Range("K113").Value = 10808.143
$ws.Range("L113").Value = 15106.214
$ws.Range("M113").Value = -8638.143
$ws.Range("N113").Value = -19446.214
$ws.Range("H126").Value = 2403.7
$ws.Range("I126").Value = 1830.25
$ws.Range("K126").Value = 5490.75
$ws.Range("M126").Value = -3020.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1930.091
$ws.Range("I9").Value = 2023.1
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 2023.1
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -1799.1
$ws.Range("N9").Value = -1448
$ws.Range("H46").Value = 4357.8423
$ws.Range("I46").Value = 1649.5
$ws.Range("J46").Value = 4676.4707
$ws.Range("K46").Value = 1649.5
$ws.Range("L46").Value = 4676.4707
$ws.Range("M46").Value = -1461.5
$ws.Range("N46").Value = -5052.4707
$ws.Range("H68").Value = 3400
$ws.Range("I68").Value = 3450
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 3450
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2701
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 3400
$ws.Range("I71").Value = 3450
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 17250
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -13506
$ws.Range("N71").Value = -22488
$ws.Range("H93").Value = 1607.5834
$ws.Range("I93").Value = 1662.8182
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1662.8182
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -414.8181999999999
$ws.Range("N93").Value = -3496

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8999
$ws.Range("J5").Value = 8999
$ws.Range("L5").Value = 8999
$ws.Range("N5").Value = -9223
$ws.Range("H132").Value = 6163.1274
$ws.Range("I132").Value = 6282.921
$ws.Range("J132").Value = 5657.3335
$ws.Range("K132").Value = 18848.763
$ws.Range("L132").Value = 16972.0005
$ws.Range("M132").Value = -16318.763
$ws.Range("N132").Value = -22032.0005
